$wb = $excel.ActiveWorkbook

# "Loan RBI, Variable Instalments" -- the Repayment Schedule sheet gains a new
# (blank) column between the existing "In Advance" (M) and "Late" (old N)
# columns, pushing "Late" -> O and "Outstanding" -> Q (with a new blank P
# between them), and it becomes the active/selected sheet (instead of Input).

$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N ("Late"); this shifts N->O,
# O->P, P->Q, matching the header/data layout of the target workbook.
$ws.Range("N1").EntireColumn.Insert()

# Make "Repayment Schedule" the active sheet/tab (was "Input").
$ws.Select() | Out-Null

# Restore the expected active cell / selection on the sheet.
$ws.Range("R8").Select() | Out-Null
